$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K5").Copy($ws.Range("L5"))
$ws.Range("K6").Copy($ws.Range("L6"))
$ws.Range("K7").Copy($ws.Range("L7"))

$r5 = $ws.Range("L5")
$r5.HorizontalAlignment = -4131
$r5.NumberFormat = """TRUE"";""TRUE"";""FALSE"""

# Now try copying format only from L5 to L6, L7
$r5.Copy($ws.Range("L6"))
$r5.Copy($ws.Range("L7"))
